$wb = $excel.ActiveWorkbook

# === 1. Reorder worksheet tabs first: Users, Recipe Types, Recipes ===
$wb.Worksheets.Item("Users").Move($wb.Worksheets.Item("Recipe Types"))

# Re-fetch fresh sheet references now that tab order has changed.
$wsUsers       = $wb.Worksheets.Item("Users")
$wsRecipeTypes = $wb.Worksheets.Item("Recipe Types")
$wsRecipes     = $wb.Worksheets.Item("Recipes ")

# === 2. Recipes sheet: mark 'temperature' as an Object (was a plain String) ===
$wsRecipes.Cells.Item(9, 2).Value = "Object"

# === 3. Users sheet: add 'email' and 'password' fields ===
$wsUsers.Cells.Item(7, 1).Value = "email"
$wsUsers.Cells.Item(7, 2).Value = "String"
$wsUsers.Cells.Item(8, 1).Value = "password"
$wsUsers.Cells.Item(8, 2).Value = "String"

# === 4. Recipes sheet: append the new 'reviews' collection fields at the bottom ===
$wsRecipes.Cells.Item(13, 1).Value = "reviews"
$wsRecipes.Cells.Item(13, 2).Value = "Array"
$wsRecipes.Cells.Item(13, 3).Value = "Object"

$wsRecipes.Cells.Item(14, 3).Value = "review_text"
$wsRecipes.Cells.Item(14, 4).Value = "String"

$wsRecipes.Cells.Item(15, 3).Value = "reviewed_by"
$wsRecipes.Cells.Item(15, 4).Value = "String"

$wsRecipes.Cells.Item(16, 3).Value = "review_date"
$wsRecipes.Cells.Item(16, 4).Value = "String"

# === 5. Recipes sheet: insert 4 rows right after 'temperature' (row 9) for its sub-fields ===
$wsRecipes.Rows.Item(10).Insert()
$wsRecipes.Rows.Item(10).Insert()
$wsRecipes.Rows.Item(10).Insert()
$wsRecipes.Rows.Item(10).Insert()

$wsRecipes.Cells.Item(10, 2).Value = "celsius"
$wsRecipes.Cells.Item(10, 3).Value = "String"

$wsRecipes.Cells.Item(11, 2).Value = "celsius_fan"
$wsRecipes.Cells.Item(11, 3).Value = "String"

$wsRecipes.Cells.Item(12, 2).Value = "fahrenheit"
$wsRecipes.Cells.Item(12, 3).Value = "String"

$wsRecipes.Cells.Item(13, 2).Value = "gas_mark"
$wsRecipes.Cells.Item(13, 3).Value = "String"

# Approximate the widened column B (now holds the nested field names)
$wsRecipes.Columns.Item(2).ColumnWidth = 11.6

# === 6. Selections on each sheet ===
$wsUsers.Range("B6:B8").Select()
$wsRecipeTypes.Range("A4").Select()
$wsRecipes.Range("C10:C13").Select()

# === 7. Active sheet / active tab ===
$wsRecipes.Activate()
